# Refresh the cryptos table's Price (D) and Volume(1h) (E) columns with the
# latest scrape. Matches commit "Updated cryptos list on Tue Oct 31
# 02:49:27 UTC 2023 with GitHub Actions".
#
# Rows 41/42 also swapped coin identity in the ranking (ARBITRUM moved
# above Aave), so those two rows get their Coin (B) and Link (C) columns
# rewritten too, in addition to Price/Volume.
#
# Every Price cell in this sheet is stored as TEXT (not a number), even
# plain-looking decimals like "228.17". A bare Range.Value assignment of
# such a string would be auto-coerced to a real number by Excel, so
# numeric-looking Price values are entered with a leading apostrophe -
# the standard "force text" prefix - which keeps the cell a string
# without leaving a visible quote character in its value/display text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '34.445.32'
$ws.Range("E2").Value = '  +0.24%  '

# Row 3
$ws.Range("D3").Value = '1.808.87'
$ws.Range("E3").Value = '  +1.20%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '''228.17'
$ws.Range("E5").Value = '  +0.78%  '

# Row 6
$ws.Range("D6").Value = '''0.580'
$ws.Range("E6").Value = '  +4.32%  '

# Row 7
$ws.Range("E7").Value = '  +0.08%  '

# Row 8
$ws.Range("D8").Value = '''35.91'
$ws.Range("E8").Value = '  +8.79%  '

# Row 9
$ws.Range("E9").Value = '  +2.95%  '

# Row 10
$ws.Range("D10").Value = '''0.0697'
$ws.Range("E10").Value = '  +1.19%  '

# Row 11
$ws.Range("E11").Value = '  +2.00%  '

# Row 12
$ws.Range("D12").Value = '2.072.10'
$ws.Range("E12").Value = '  +1.23%  '

# Row 13
$ws.Range("D13").Value = '''11.52'
$ws.Range("E13").Value = '  +3.15%  '

# Row 14
$ws.Range("D14").Value = '1.825.06'
$ws.Range("E14").Value = '  +2.08%  '

# Row 15
$ws.Range("E15").Value = '  +2.06%  '

# Row 16
$ws.Range("D16").Value = '''4.53'
$ws.Range("E16").Value = '  +5.62%  '

# Row 17
$ws.Range("D17").Value = '34.458.02'
$ws.Range("E17").Value = '  +0.29%  '

# Row 18
$ws.Range("D18").Value = '''69.32'
$ws.Range("E18").Value = '  +1.11%  '

# Row 19
$ws.Range("D19").Value = '''247.05'
$ws.Range("E19").Value = '  +0.69%  '

# Row 20
$ws.Range("E20").Value = '  +0.45%  '

# Row 21
$ws.Range("D21").Value = '''11.56'
$ws.Range("E21").Value = '  +2.43%  '

# Row 22
$ws.Range("E22").Value = '  +0.06%  '

# Row 23
$ws.Range("E23").Value = '  +1.41%  '

# Row 24
$ws.Range("D24").Value = '''171.76'
$ws.Range("E24").Value = '  +1.58%  '

# Row 25
$ws.Range("E25").Value = '  +3.23%  '

# Row 26
$ws.Range("D26").Value = '''7.99'
$ws.Range("E26").Value = '  +8.84%  '

# Row 27
$ws.Range("D27").Value = '''17.05'
$ws.Range("E27").Value = '  +3.14%  '

# Row 28
$ws.Range("E28").Value = '  +3.79%  '

# Row 29
$ws.Range("E29").Value = '  -0.05%  '

# Row 30
$ws.Range("D30").Value = '''4.07'
$ws.Range("E30").Value = '  +1.06%  '

# Row 31
$ws.Range("D31").Value = '''0.0534'
$ws.Range("E31").Value = '  +1.75%  '

# Row 32
$ws.Range("E32").Value = '  +2.27%  '

# Row 33
$ws.Range("E33").Value = '  +0.56%  '

# Row 34
$ws.Range("E34").Value = '  +1.58%  '

# Row 35
$ws.Range("D35").Value = '1.401.02'
$ws.Range("E35").Value = '  -0.84%  '

# Row 36
$ws.Range("D36").Value = '''0.678'
$ws.Range("E36").Value = '  -0.66%  '

# Row 37
$ws.Range("D37").Value = '''2.50'
$ws.Range("E37").Value = '  -2.91%  '

# Row 38
$ws.Range("E38").Value = '  +0.50%  '

# Row 39
$ws.Range("E39").Value = '  +0.78%  '

# Row 40
$ws.Range("D40").Value = '''1.25'
$ws.Range("E40").Value = '  +13.02%  '

# Row 41
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '''0.968'
$ws.Range("E41").Value = '  +2.87%  '

# Row 42
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '''83.08'
$ws.Range("E42").Value = '  -1.34%  '

# Row 43
$ws.Range("E43").Value = '  +1.50%  '

# Row 44
$ws.Range("E44").Value = '  +0.13%  '

# Row 45
$ws.Range("D45").Value = '''13.41'
$ws.Range("E45").Value = '  -4.10%  '

# Row 46
$ws.Range("E46").Value = '  -1.18%  '

# Row 47
$ws.Range("D47").Value = '''0.0506'
$ws.Range("E47").Value = '  -3.96%  '

# Row 48
$ws.Range("D48").Value = '1.970.98'
$ws.Range("E48").Value = '  +1.12%  '

# Row 49
$ws.Range("D49").Value = '''105.73'
$ws.Range("E49").Value = '  +0.45%  '

# Row 50
$ws.Range("E50").Value = '  +0.11%  '

# Row 51
$ws.Range("E51").Value = '  +1.79%  '
